$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new record at row 139 ---
# This pushes the existing rows 139..211 down to 140..212.
$ws.Rows.Item(139).Insert()

$ws.Range("A139").Value = 3
$ws.Range("B139").Value = "Femacal de La Calera"
$ws.Range("C139").Value = "Coquimbo"
$ws.Range("D139").Value = 44846
$ws.Range("E139").Value = 5
$ws.Range("F139").Value = 100112026
$ws.Range("G139").Value = "Haba"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 85
$ws.Range("K139").Value = 8000
$ws.Range("L139").Value = 8500
$ws.Range("M139").Value = 8265
$ws.Range("N139").Value = "$/saco 25 kilos"
$ws.Range("O139").Value = "Provincia de Limarí"
$ws.Range("P139").Value = 331
$ws.Range("Q139").Value = 25
$ws.Range("R139").Value = "Hortaliza"

# --- Insert second new record at row 147 ---
# At this point the original row 146 (before the first insert) now sits at
# row 147, so inserting here pushes it (and everything after) down by one
# more row, landing the new record at row 147 and the remaining original
# data two rows below its original position overall.
$ws.Rows.Item(147).Insert()

$ws.Range("A147").Value = 3
$ws.Range("B147").Value = "Femacal de La Calera"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = 44845
$ws.Range("E147").Value = 5
$ws.Range("F147").Value = 100112026
$ws.Range("G147").Value = "Haba"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 145
$ws.Range("K147").Value = 8500
$ws.Range("L147").Value = 9000
$ws.Range("M147").Value = 8759
$ws.Range("N147").Value = "$/saco 25 kilos"
$ws.Range("O147").Value = "Provincia de Limarí"
$ws.Range("P147").Value = 350
$ws.Range("Q147").Value = 25
$ws.Range("R147").Value = "Hortaliza"
